# Refresh the cryptos price list (rows 2-51, columns B-E) to the
# latest snapshot values from the commit diff.
#
# Some new "Price" values are plain decimal numbers (e.g. "0.606").
# Excel's Range.Value setter auto-detects such text as a number, which
# would change the cell's stored type away from the original text/
# inline-string representation. To keep these as text we enter them
# with a leading apostrophe (Excel's "treat as text" entry prefix,
# which is not stored as part of the value) and then reset the cell's
# style back to "Normal" so no stray number-format/style is left on
# the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.640.11'
$ws.Range('E2').Value = '  +3.67%  '
$ws.Range('D3').Value = '2.651.97'
$ws.Range('E3').Value = '  +1.19%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''571.43'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +7.15%  '
$ws.Range('D6').Value = '''147.31'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.44%  '
$ws.Range('E7').Value = '  -0.38%  '
$ws.Range('D8').Value = '''0.606'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.64%  '
$ws.Range('E9').Value = '  -0.66%  '
$ws.Range('E10').Value = '  +4.53%  '
$ws.Range('D11').Value = '''0.144'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.54%  '
$ws.Range('D12').Value = '''0.344'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.20%  '
$ws.Range('D13').Value = '3.124.95'
$ws.Range('E13').Value = '  +1.26%  '
$ws.Range('D14').Value = '60.659.73'
$ws.Range('E14').Value = '  +3.84%  '
$ws.Range('D15').Value = '''21.84'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.85%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '''0.0000138'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.91%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.672.11'
$ws.Range('E17').Value = '  +1.40%  '
$ws.Range('D18').Value = '''4.55'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.82%  '
$ws.Range('D19').Value = '''344.64'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.13%  '
$ws.Range('D20').Value = '''10.47'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.50%  '
$ws.Range('D21').Value = '''6.43'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.37%  '
$ws.Range('D22').Value = '''5.83'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.04%  '
$ws.Range('D23').Value = '''0.999'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').Value = '''66.70'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.47%  '
$ws.Range('E25').Value = '  +7.25%  '
$ws.Range('E26').Value = '  +1.67%  '
$ws.Range('D27').Value = '''0.991'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.87%  '
$ws.Range('E28').Value = '  +4.57%  '
$ws.Range('D29').Value = '0.0₃0791'
$ws.Range('E29').Value = '  +8.01%  '
$ws.Range('E30').Value = '  -0.10%  '
$ws.Range('E31').Value = '  +5.23%  '
$ws.Range('D32').Value = '''6.14'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.13%  '
$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').Value = '''155.31'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.41%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = '''19.31'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.01%  '
$ws.Range('E35').Value = '  +6.31%  '
$ws.Range('D36').Value = '''0.918'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.61%  '
$ws.Range('E37').Value = '  +8.19%  '
$ws.Range('D38').Value = '''0.917'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +13.48%  '
$ws.Range('D39').Value = '''37.65'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.30%  '
$ws.Range('E40').Value = '  +7.92%  '
$ws.Range('D41').Value = '''309.27'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +10.61%  '
$ws.Range('E42').Value = '  +3.20%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = '''0.993'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.62%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').Value = '''0.609'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.74%  '
$ws.Range('D45').Value = '''0.0981'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.94%  '
$ws.Range('D46').Value = '''0.0550'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.43%  '
$ws.Range('D47').Value = '''19.49'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.83%  '
$ws.Range('D48').Value = '''10.66'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.11%  '
$ws.Range('D49').Value = '''126.39'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +11.98%  '
$ws.Range('E50').Value = '  +5.34%  '
$ws.Range('D51').Value = '''4.76'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.63%  '
